# Actualización automática 2025-08-25 12:40:09
#
# This script applies the diff to both worksheets:
#  - "VENTAS POR GRUPO" (sheet 1, columns A:R)
#  - "VENTA MENSUAL"    (sheet 2, columns A:G)
#
# In both sheets, the roster of (ASESOR, CLIENTE) rows grows from 3 to 6
# clients for "ILLER LOPEZ ROBERTO FERNANDO":
#   old row2 CLIENTE "COELLO TRONCOSO JOSE GREGORIO"       -> "ASES GAVILANEZ FAUSTO HERNAN"
#   old row3 CLIENTE "COMERCIAL LUNA PAZMIÑO CIA. LTDA."   -> "BRAVO MONTENEGRO DANIEL ANDRES"
#   three brand-new rows are inserted before the old row4, for clients:
#       COELLO TRONCOSO JOSE GREGORIO
#       COMERCIAL LUNA PAZMIÑO CIA. LTDA.
#       LATACELA ZUÑIGA JUAN FERNANDO
#   old row4 (VIEJO RIVAS MAYRA ANABELLE) shifts down to row7
#   the trailing totals row shifts from row5 to row8, and its "0 de 3"
#   label (sheet1) becomes "0 de 6"; sheet2 totals pick up the new
#   LATACELA ZUÑIGA JUAN FERNANDO figures (171.19 / 0 / 0 / 0 / 1500).

$wb = $excel.ActiveWorkbook

$newClients = @(
    "COELLO TRONCOSO JOSE GREGORIO",
    "COMERCIAL LUNA PAZMIÑO CIA. LTDA.",
    "LATACELA ZUÑIGA JUAN FERNANDO"
)
$advisor = "ILLER LOPEZ ROBERTO FERNANDO"

# ---------------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"  (data columns C..R)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Rename two existing clients.
$ws1.Range("B2").Value = "ASES GAVILANEZ FAUSTO HERNAN"
$ws1.Range("B3").Value = "BRAVO MONTENEGRO DANIEL ANDRES"

# Insert 3 new blank rows before row 4 (pushes old row4 -> row7, old row5 -> row8).
$ws1.Range("A4:R6").Insert()

# Fill the three new rows with advisor/client names and zeroed metrics.
for ($i = 0; $i -lt $newClients.Length; $i++) {
    $r = 4 + $i
    $ws1.Range("A$r").Value = $advisor
    $ws1.Range("B$r").Value = $newClients[$i]
    for ($c = 3; $c -le 18; $c++) {
        $ws1.Cells.Item($r, $c).Value = 0
    }
}

# Update the trailing "0 de 3" -> "0 de 6" summary row (now row 8).
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(8, $c).Value = "0 de 6"
}

# ---------------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"  (data columns C..G)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Widen column C slightly (displayed width 10 -> 12).
$ws2.Columns.Item(3).ColumnWidth = 11.14

# Rename two existing clients.
$ws2.Range("B2").Value = "ASES GAVILANEZ FAUSTO HERNAN"
$ws2.Range("B3").Value = "BRAVO MONTENEGRO DANIEL ANDRES"

# Insert 3 new blank rows before row 4 (pushes old row4 -> row7, old row5 -> row8).
$ws2.Range("A4:G6").Insert()

# Fill the three new rows with advisor/client names and zeroed metrics.
for ($i = 0; $i -lt $newClients.Length; $i++) {
    $r = 4 + $i
    $ws2.Range("A$r").Value = $advisor
    $ws2.Range("B$r").Value = $newClients[$i]
    for ($c = 3; $c -le 7; $c++) {
        $ws2.Cells.Item($r, $c).Value = 0
    }
}

# LATACELA ZUÑIGA JUAN FERNANDO (row 6) carries real sales figures.
$ws2.Range("C6").Value = 171.19
$ws2.Range("D6").Value = 0
$ws2.Range("E6").Value = 0
$ws2.Range("F6").Value = 0
$ws2.Range("G6").Value = 1500

# Update the trailing totals row (now row 8) to reflect the new figures.
$ws2.Range("C8").Value = 171.19
$ws2.Range("D8").Value = 0
$ws2.Range("E8").Value = 0
$ws2.Range("F8").Value = 0
$ws2.Range("G8").Value = 1500
